# Auto-generated edit script: update crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.872.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.24%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.357.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'503.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'130.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.13%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.19%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.45%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.372.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0965"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.89%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.52%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.02%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.321"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.777.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.96%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'55.800.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.04%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.372.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.34%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'9.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'310.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.23%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.46%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.43%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.28%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.12%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.22%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.04%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'171.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.93%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0710"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.12%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.78%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.47%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.15%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.843"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -4.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'36.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.85%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'PolygonEcosystemToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.370"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Filecoin"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.24%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'RenderToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'4.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.80%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Aave"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'125.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.79%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Mantle"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.557"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.93%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Stellar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0895"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.59%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Bittensor"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'240.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Hedera"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0478"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.34%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'16.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.11%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'VeChain"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.96%  "
$ws.Range("E51").Style = "Normal"
